$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (kr-vs-kp): update the NOTE text
$ws.Range("K4").Value = "necessarie tutte le combinazioni"

# Row 6 (Balloons): replace the stray leftover note text with the real note,
# and give it the same formatting as the other NOTE cells (left aligned)
$ws.Range("K6").Value = "necessarie tutte le combinazioni"
$ws.Range("K6").HorizontalAlignment = -4131  # xlLeft

# Row 8 (Sepsis): clear the NOTE text
$ws.Range("K8").ClearContents()

# Re-fit column K width to the new (longer) text
$ws.Columns.Item(11).ColumnWidth = 26.66
